$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '62.313.01'
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.889.38'
$ws.Range('E3').Value = '  -0.64%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '565.88'
$ws.Range('E5').Value = '  -3.53%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.18'
$ws.Range('E6').Value = '  -2.19%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -1.41%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.887.74'
$ws.Range('E9').Value = '  -0.67%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.86'
$ws.Range('E10').Value = '  +0.48%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.148'
$ws.Range('E11').Value = '  -1.28%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.429'
$ws.Range('E12').Value = '  -1.14%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000236'
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '32.64'
$ws.Range('E14').Value = '  -0.30%  '
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.376.45'
$ws.Range('E16').Value = '  -0.47%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '62.362.01'
$ws.Range('E17').Value = '  +0.82%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.56'
$ws.Range('E18').Value = '  -1.20%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '2.895.83'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '425.79'
$ws.Range('E20').Value = '  -2.35%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.98'
$ws.Range('E21').Value = '  -2.83%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.655'
$ws.Range('E22').Value = '  -0.62%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.84'
$ws.Range('E23').Value = '  -1.64%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '78.44'
$ws.Range('E24').Value = '  -2.66%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '11.81'
$ws.Range('E25').Value = '  -1.26%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.99'
$ws.Range('E27').Value = '  -2.19%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.00'
$ws.Range('E28').Value = '  -3.34%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0000108'
$ws.Range('E29').Value = '  +1.42%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.06'
$ws.Range('E30').Value = '  -1.07%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.47'
$ws.Range('E31').Value = '  -3.40%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.00'
$ws.Range('E32').Value = '  -4.86%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '25.62'
$ws.Range('E34').Value = '  -0.97%  '
$ws.Range('E35').Value = '  -3.46%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.943'
$ws.Range('E36').Value = '  -3.23%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.36'
$ws.Range('E37').Value = '  -2.60%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '48.42'
$ws.Range('E38').Value = '  -1.47%  '
$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.89'
$ws.Range('E39').Value = '  -4.59%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.88'
$ws.Range('E40').Value = '  -5.59%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '41.11'
$ws.Range('E41').Value = '  +5.62%  '
$ws.Range('E42').Value = '  -2.07%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.00'
$ws.Range('E43').Value = '  -4.60%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.713.31'
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.265'
$ws.Range('E45').Value = '  -2.76%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '133.27'
$ws.Range('E46').Value = '  -0.54%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0336'
$ws.Range('E47').Value = '  +0.27%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '353.24'
$ws.Range('E48').Value = '  +3.30%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.000215'
$ws.Range('E50').Value = '  +11.10%  '
$ws.Range('E51').Value = '  -0.96%  '
